$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Date column (D2:D5) to the new date serial value 45852
$ws.Range("D2").Value = 45852
$ws.Range("D3").Value = 45852
$ws.Range("D4").Value = 45852
$ws.Range("D5").Value = 45852

# Update the Time of Test column (E2:E5) as text strings
$ws.Range("E2").Value = "01:35:30"
$ws.Range("E3").Value = "01:35:37"
$ws.Range("E4").Value = "01:35:37"
$ws.Range("E5").Value = "01:35:44"
